# TC_149 test data update - "Updated test data as per new implemenation"
#
# 1. Rename the column headers on the "Add Panels" sheet (row 8, cols N:O)
#    from "Battery Alarm (A)" / "Battery Standby (A)"
#    to   "Alarm Current(A)" / "Standby Current(A)"
# 2. Move the sheet's active selection from N7:O8 down to N8:O8
#    (matching the new location of the renamed values).
# 3. Re-save the workbook against the corrected working-copy path
#    (C:\Work\... -> C:\work\...) to mirror the updated absPath metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# --- 1. Update the relabeled header text -------------------------------
$ws.Range("N8").Value = "Alarm Current(A)"
$ws.Range("O8").Value = "Standby Current(A)"

# --- 2. Update the active selection on the sheet ------------------------
$ws.Activate()
$ws.Range("N8:O8").Select()

# --- 3. Reflect the corrected case of the working directory -------------
try {
    $wb.SaveAs("C:\work\consys-uiauto\Test Data\TC_149_Verify_Battery_Standby_And_Alarm_Load_On_Addition_Deletion_Of_Accessories.xlsx")
} catch {
    # Best effort only - some hosts do not track the originating path.
}
